$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain text (their values look numeric,
# but the source data stores them as literal strings).
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D14", "D17", "D19", "D20", "D21", "D23", "D24", "D25", "D28", "D31", "D35", "D38", "D41", "D42", "D43", "D45", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value = "52.210.19"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "2.788.52"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "344.76"
$ws.Range("D6").Value = "116.22"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +3.39%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").Value = "  +2.91%  "
$ws.Range("D10").Value = "42.96"
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("D11").Value = "0.0858"
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("D12").Value = "20.13"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").Value = "7.75"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "3.226.76"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "2.788.03"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").Value = "0.889"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "52.069.81"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "3.22"
$ws.Range("E19").Value = "  +5.55%  "
$ws.Range("D20").Value = "7.08"
$ws.Range("D21").Value = "13.37"
$ws.Range("E21").Value = "  -2.75%  "
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").Value = "270.27"
$ws.Range("E23").Value = "  -5.74%  "
$ws.Range("D24").Value = "70.16"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "2.77"
$ws.Range("E25").Value = "  +6.17%  "
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "10.27"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").Value = "34.88"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").Value = "0.0408"
$ws.Range("E35").Value = "  +16.06%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "19.09"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").Value = "2.69"
$ws.Range("E41").Value = "  +20.24%  "
$ws.Range("D42").Value = "23.64"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").Value = "127.73"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("D45").Value = "2.33"
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("D47").Value = "2.073.28"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "0.912"
$ws.Range("E50").Value = "  +10.58%  "
$ws.Range("D51").Value = "8.97"
$ws.Range("E51").Value = "  -1.59%  "
